# Regenerate the handback-status report: the two tracked files were
# re-handed-back, producing new GUID-named files and fresh xliff
# handoff/handback timestamps.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet --------------------------------------------------
$wsOverview.Range("A2").Value = "67f62cdd-6a9a-4960-9101-b642c3f49ce7.md"
$wsOverview.Range("B2").Value = "e2e\67f62cdd-6a9a-4960-9101-b642c3f49ce7.md"
$wsOverview.Range("G2").Value = "2016-08-17 10:59:20"

$wsOverview.Range("A3").Value = "ffff50ea4798-9698-448d-9600-0478c9c41ecf.md"
$wsOverview.Range("B3").Value = "e2e\ffff50ea4798-9698-448d-9600-0478c9c41ecf.md"
$wsOverview.Range("G3").Value = "2016-08-17 10:59:20"

# --- zh-cn sheet -------------------------------------------------------
$wsZhCn.Range("A2").Value = "67f62cdd-6a9a-4960-9101-b642c3f49ce7.md"
$wsZhCn.Range("G2").Value = "67f62cdd-6a9a-4960-9101-b642c3f49ce7.50748a98f39a2216731aa09752af2d6c3591594f.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-17 10:59:16"
$wsZhCn.Range("I2").Value = "67f62cdd-6a9a-4960-9101-b642c3f49ce7.md"
$wsZhCn.Range("J2").Value = "67f62cdd-6a9a-4960-9101-b642c3f49ce7.50748a98f39a2216731aa09752af2d6c3591594f.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-17 10:59:33"

$wsZhCn.Range("A3").Value = "ffff50ea4798-9698-448d-9600-0478c9c41ecf.md"
$wsZhCn.Range("G3").Value = "67f62cdd-6a9a-4960-9101-b642c3f49ce7.50748a98f39a2216731aa09752af2d6c3591594f.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-17 10:59:16"
$wsZhCn.Range("I3").Value = "ffff50ea4798-9698-448d-9600-0478c9c41ecf.md"
$wsZhCn.Range("J3").Value = "67f62cdd-6a9a-4960-9101-b642c3f49ce7.50748a98f39a2216731aa09752af2d6c3591594f.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-17 10:59:33"

# --- de-de sheet -------------------------------------------------------
$wsDeDe.Range("A2").Value = "67f62cdd-6a9a-4960-9101-b642c3f49ce7.md"
$wsDeDe.Range("G2").Value = "67f62cdd-6a9a-4960-9101-b642c3f49ce7.50748a98f39a2216731aa09752af2d6c3591594f.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-17 10:59:20"
$wsDeDe.Range("I2").Value = "67f62cdd-6a9a-4960-9101-b642c3f49ce7.md"
$wsDeDe.Range("J2").Value = "67f62cdd-6a9a-4960-9101-b642c3f49ce7.50748a98f39a2216731aa09752af2d6c3591594f.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-17 10:59:41"

$wsDeDe.Range("A3").Value = "ffff50ea4798-9698-448d-9600-0478c9c41ecf.md"
$wsDeDe.Range("G3").Value = "67f62cdd-6a9a-4960-9101-b642c3f49ce7.50748a98f39a2216731aa09752af2d6c3591594f.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-17 10:59:20"
$wsDeDe.Range("I3").Value = "ffff50ea4798-9698-448d-9600-0478c9c41ecf.md"
$wsDeDe.Range("J3").Value = "67f62cdd-6a9a-4960-9101-b642c3f49ce7.50748a98f39a2216731aa09752af2d6c3591594f.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-17 10:59:41"

# --- Hyperlink display text updates ------------------------------------
# The hyperlink targets (r:id) are untouched; only the visible display
# text (which mirrors the renamed files) needs to change.
foreach ($ws in @($wsOverview, $wsZhCn, $wsDeDe)) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.TextToDisplay -like "*871c2a8f-e428-4c32-a3d7-4f079e4772ac.md") {
            $hl.TextToDisplay = $hl.TextToDisplay.Replace("871c2a8f-e428-4c32-a3d7-4f079e4772ac.md", "67f62cdd-6a9a-4960-9101-b642c3f49ce7.md")
        } elseif ($hl.TextToDisplay -like "*a278e830-8c0f-430d-b262-e3153c0360f6.md") {
            $hl.TextToDisplay = $hl.TextToDisplay.Replace("a278e830-8c0f-430d-b262-e3153c0360f6.md", "ffff50ea4798-9698-448d-9600-0478c9c41ecf.md")
        }
    }
}
